$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 46, shifting rows 46-71 down to 47-72.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new record.
$ws.Cells.Item(46, 1).Value = 5
$ws.Cells.Item(46, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(46, 3).Value = "Maule"
$ws.Cells.Item(46, 4).Value = 44529
$ws.Cells.Item(46, 5).Value = 7
$ws.Cells.Item(46, 6).Value = 100112022
$ws.Cells.Item(46, 7).Value = "Arveja Verde"
$ws.Cells.Item(46, 8).Value = "Perfection"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 400
$ws.Cells.Item(46, 11).Value = 13000
$ws.Cells.Item(46, 12).Value = 13000
$ws.Cells.Item(46, 13).Value = 13000
$ws.Cells.Item(46, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(46, 15).Value = "Región del Maule"
$ws.Cells.Item(46, 16).Value = 520
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = "Hortaliza"
